# This script reproduces the "Common: Added some initial data for vendors,
# atomizers, mods" commit: 5 new vendor names are added to the "vendors"
# sheet, and the (already-sorted) column is kept alphabetically sorted, so
# most existing rows shift down by one or more positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vendors")

# --- Pass 1: re-seat every pre-existing vendor name into its new (post-sort)
# row. These strings already exist in the workbook, so re-assigning them does
# not create new shared-string entries or change any existing index/order.
$existingVendors = @{
    2 = "3 Baccos"
    3 = "Acrohm"
    5 = "Adams Vape"
    6 = "Al Carlo"
    7 = "Ambition MODS"
    8 = "Aramax"
    9 = "Asmodus"
    10 = "Aspire"
    11 = "Aviator MODS"
    12 = "BP MODS"
    13 = "CoilArt"
    14 = "Coilology"
    15 = "Cthulhu"
    16 = "Dekang"
    17 = "Diamond Mist"
    18 = "Dicodes"
    19 = "Digiflavor"
    20 = "Dinner Lady"
    21 = "Dotmod"
    22 = "Dovpo"
    23 = "Ecoliquid"
    24 = "Ehpro"
    25 = "ElcigART"
    26 = "Emporio"
    27 = "Fakirsmods"
    28 = "Flavormonks"
    29 = "GeekVape"
    31 = "Imperia"
    32 = "Infamous"
    33 = "iSmoka-Eleaf"
    34 = "IVG"
    35 = "Joyetech"
    36 = "Lost Vape"
    37 = "Mark Bugs"
    38 = "Mechlyfe"
    39 = "My Vape"
    40 = "Ohm Boy"
    41 = "PJ Empire"
    42 = "qp Design"
    43 = "Rev Tech"
    44 = "SMArt Mods"
    45 = "SmokerStore"
    46 = "Smoktech"
    47 = "Squape"
    48 = "Squid industries"
    49 = "Sunbox"
    51 = "SXK"
    52 = "The Crazy Wire"
    53 = "Ultroner"
    55 = "UWELL"
    56 = "Vandy Vape"
    57 = "Vap Extreme"
    58 = "Vapefly"
    59 = "Vapor Giant"
    60 = "Vaptio"
    61 = "VGOD"
    63 = "VooPoo"
    64 = "Wotofo"
    65 = "Yihi"
}
foreach ($row in $existingVendors.Keys) {
    $ws.Cells.Item($row, 1).Value = $existingVendors[$row]
}

# --- Pass 2: write the five brand-new vendor names, in the order they were
# first entered into the sheet (preserves the expected shared-string order).
$ws.Cells.Item(4, 1).Value = "Acrossvape"
$ws.Cells.Item(30, 1).Value = "HussarVape"
$ws.Cells.Item(62, 1).Value = "Vicious Ant"
$ws.Cells.Item(50, 1).Value = "SvoëMesto"
$ws.Cells.Item(54, 1).Value = "Unknown"

# --- Tidy up view state to match the edited region (best effort).
$ws.Range("A55").Select()

try {
    $sort = $ws.Sort
    $sort.SortFields.Clear()
    $sort.SortFields.Add($ws.Range("A2:A65"))
    $sort.SetRange($ws.Range("A1:A65"))
    $sort.Header = 1
    $sort.Apply()
    Write-Host "sort applied"
} catch {
    Write-Host ("sort apply failed: " + $_)
}
